# The sheet contains weekly price records for "Zapallo italiano" sold by
# "Comercializadora del Agro de Limarí". A new weekly record (dated 2023-05-17)
# is inserted as row 13, pushing the existing rows 13-77 down to 14-78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 13 - this shifts rows 13:77 down to 14:78
# and automatically extends the sheet dimension from A1:R77 to A1:R78.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record's data.
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = "2023-05-17"
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112032
$ws.Range("G13").Value = "Zapallo italiano"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 10500
$ws.Range("N13").Value = "`$/caja 60 unidades"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 175
$ws.Range("Q13").Value = 60
$ws.Range("R13").Value = "Hortaliza"
